$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 4 new rows right after the "Silniki krokowe" section header row
#    (row 3), growing that section from 3 rows to 7 rows. Excel inherits the
#    formatting of the row above for the newly inserted rows.
# ---------------------------------------------------------------------------
$ws.Rows("4:7").Insert()
$ws.Rows("4:7").RowHeight = 36

# The old D4 note ("nie bardzo wiadomo...") shifted down to D8 along with its
# row - it belongs with "Przelozenie zebatek", which is being moved to the
# bottom of the section (row 9), so clear it from its shifted location.
$ws.Range("D8").Clear()

# ---------------------------------------------------------------------------
# 2. Rewrite the "Silniki krokowe" section content to match the final order:
#    new ESP/Halla related tasks inserted, "Przelozenie zebatek" (+ its note)
#    moved to become the last item of the section.
# ---------------------------------------------------------------------------
$ws.Range("C4").Value = "Potestować konfigurację z czujnikiem Halla (coś w stylu PID)"
$ws.Range("C5").Value = "Poeksperymentować z trybami pracy driverów w silnikach"
$ws.Range("C6").Value = "Podłączyć scp + raspberry przez Serial"
$ws.Range("C7").Value = "Przelutować pajęczaka na płytkę"
$ws.Range("C8").Value = "Przeprojektować model w inventorze, by obsługiwał silniki i Halla"
$ws.Range("C9").Value = "Przełożenie zębatek wymaga aktualizacji"
$ws.Range("D9").Value = "nie bardzo wiadomo jakiego przełożenia użyć - do sprawdzenia z realnym obciążeniem"

# D9 keeps the same look the note had before (left / center, wrapped).
$ws.Range("D9").HorizontalAlignment = -4131
$ws.Range("D9").VerticalAlignment = -4108
$ws.Range("D9").WrapText = $true

# ---------------------------------------------------------------------------
# 3. Highlight (green) the "pinned" rows: the first item under "Silniki
#    krokowe" (C3), the whole "Lutowanie" header row (B10:C10), and refresh
#    the existing green highlight on the "Audio" items (C11:C12) to the new
#    solid RGB green instead of the old theme-based green.
# ---------------------------------------------------------------------------
$green = 5296274  # RGB(146, 208, 80) / FF92D050

$ws.Range("C3").Interior.Color = $green
$ws.Range("C3").VerticalAlignment = -4108
$ws.Range("C3").WrapText = $true

$ws.Range("B10").Interior.Color = $green
$ws.Range("B10").HorizontalAlignment = -4152
$ws.Range("B10").VerticalAlignment = -4108

$ws.Range("C10").Interior.Color = $green
$ws.Range("C10").VerticalAlignment = -4108
$ws.Range("C10").WrapText = $true

$ws.Range("C11:C12").Interior.Color = $green
$ws.Range("C11:C12").HorizontalAlignment = -4131
$ws.Range("C11:C12").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 4. Misc bookkeeping to match the authored workbook exactly.
# ---------------------------------------------------------------------------
$ws.Range("D6").Select()
